$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '29.375.68'
$ws.Range("D2").ClearFormats()
$ws.Range("E2").Value = '  +0.60%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.874.47'
$ws.Range("D3").ClearFormats()
$ws.Range("E3").Value = '  +0.72%  '
$ws.Range("E4").Value = '  +0.07%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '0.7115'
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = '  -0.26%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '242.04'
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = '  +0.76%  '
$ws.Range("E7").Value = '  -0.13%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3114'
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = '  +1.06%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.07778'
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = '  +1.26%  '
$ws.Range("E10").Value = '  +0.10%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.08465'
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = '  +1.75%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '1.866.37'
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = '  -2.46%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '5.241'
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = '  +0.61%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.7121'
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = '  -0.54%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '91.17'
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = '  +0.54%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '29.377.16'
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = '  +0.21%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.000008249'
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = '  +5.78%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '6.038'
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = '  +1.41%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '240.86'
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = '  -0.82%  '
$ws.Range("E20").Value = '  +0.86%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '2.122.02'
$ws.Range("D21").ClearFormats()
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '1.000'
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = '  -0.07%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '7.795'
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = '  -2.30%  '
$ws.Range("E24").Value = '  +0.00%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '0.1611'
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = '  -0.13%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '163.91'
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = '  +0.70%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '9.066'
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = '  +1.92%  '
$ws.Range("E28").Value = '  -0.48%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '1.512'
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = '  +1.14%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '4.427'
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = '  -0.18%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '4.314'
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = '  +1.66%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '1.283'
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = '  -4.30%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.05298'
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = '  +2.17%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.938'
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = '  +0.80%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.178'
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = '  +0.65%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.7472'
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = '  -5.32%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '2.694'
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = '  +0.38%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.01870'
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = '  +0.83%  '
$ws.Range("B39").Value = 'MXToken'
$ws.Range("C39").Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '2.722'
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = '  +1.18%  '
$ws.Range("B40").Value = 'Maker'
$ws.Range("C40").Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '1.205.39'
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = '  +2.21%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '6.450'
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = '  +3.36%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.8881'
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = '  -1.58%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '72.81'
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = '  +0.04%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '108.98'
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = '  +6.49%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '1.000'
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = '  +0.03%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '2.020.50'
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = '  -2.38%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '1.818'
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = '  +2.45%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.5210'
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = '  +0.20%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.00000000123'
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = '  +4.70%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '9.386'
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = '  +0.48%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.4318'
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = '  +1.00%  '
